# Generate Report for Handback
# Update timestamps / priority values produced by a new handback report run.

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-09 01:39:24"
$wsOverview.Range("G3").Value = "2016-11-09 01:39:24"

# --- "zh-cn" sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-11-09 01:39:08"
$wsZhCn.Range("H3").Value = "2016-11-09 01:39:08"
$wsZhCn.Range("K2").Value = "2016-11-09 01:40:01"
$wsZhCn.Range("K3").Value = "2016-11-09 01:40:01"

# --- "de-de" sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-09 01:39:24"
$wsDeDe.Range("H3").Value = "2016-11-09 01:39:24"
$wsDeDe.Range("K2").Value = "2016-11-09 01:40:21"
$wsDeDe.Range("K3").Value = "2016-11-09 01:40:21"
